$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting existing rows 14:115 down to 15:116
$ws.Rows.Item(14).Insert()

# Copy the unchanged "metadata" columns from the row that used to be row 14
# (now shifted to row 15) into the new row 14.
$ws.Cells.Item(14, 1).Value = $ws.Cells.Item(15, 1).Value()    # A Mercado ID
$ws.Cells.Item(14, 2).Value = $ws.Cells.Item(15, 2).Value()    # B Mercado
$ws.Cells.Item(14, 3).Value = $ws.Cells.Item(15, 3).Value()    # C Region
$ws.Cells.Item(14, 5).Value = $ws.Cells.Item(15, 5).Value()    # E Codreg
$ws.Cells.Item(14, 6).Value = $ws.Cells.Item(15, 6).Value()    # F Categoria ID
$ws.Cells.Item(14, 7).Value = $ws.Cells.Item(15, 7).Value()    # G Categoria
$ws.Cells.Item(14, 8).Value = $ws.Cells.Item(15, 8).Value()    # H Variedad
$ws.Cells.Item(14, 9).Value = $ws.Cells.Item(15, 9).Value()    # I Calidad
$ws.Cells.Item(14, 14).Value = $ws.Cells.Item(15, 14).Value()  # N Unidad de comercializacion
$ws.Cells.Item(14, 15).Value = $ws.Cells.Item(15, 15).Value()  # O Origen
$ws.Cells.Item(14, 17).Value = $ws.Cells.Item(15, 17).Value()  # Q Kg o Unidades
$ws.Cells.Item(14, 18).Value = $ws.Cells.Item(15, 18).Value()  # R Clasificacion

# Set the new values for the newly inserted row 14
$ws.Cells.Item(14, 4).Value = 45168    # D Fecha
$ws.Cells.Item(14, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"   # keep column D date formatting
$ws.Cells.Item(14, 10).Value = 16      # J Volumen
$ws.Cells.Item(14, 11).Value = 21000   # K Precio minimo
$ws.Cells.Item(14, 12).Value = 21000   # L Precio maximo
$ws.Cells.Item(14, 13).Value = 21000   # M Precio promedio ponderado
$ws.Cells.Item(14, 16).Value = 7000    # P Precio $/Kg
